$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update room-id (A) / service-id (B) values: prefix with hotel 1 id ---
$ws.Range("A2").Value = 10001101
$ws.Range("B2").Value = 1000110101

$ws.Range("A3").Value = 10001102
$ws.Range("B3").Value = 1000110102

$ws.Range("A4").Value = 10001103
$ws.Range("B4").Value = 1000110103

$ws.Range("A5").Value = 10001104
$ws.Range("B5").Value = 1000110104

$ws.Range("A6").Value = 10001105
$ws.Range("B6").Value = 1000110105

$ws.Range("A7").Value = 10001106
$ws.Range("B7").Value = 1000110106

$ws.Range("A8").Value = 10001107
$ws.Range("B8").Value = 1000110107

$ws.Range("A9").Value = 10001108
$ws.Range("B9").Value = 1000110108

$ws.Range("A10").Value = 10001109
$ws.Range("B10").Value = 1000110109

$ws.Range("A11").Value = 10001110
$ws.Range("B11").Value = 1000110110

# --- Column widths now needed to fit the longer ids (bestFit-style autosize) ---
$ws.Columns.Item(1).ColumnWidth = 11.3
$ws.Columns.Item(2).ColumnWidth = 10.3

# --- Selection moves to B11 ---
$ws.Range("B11").Select() | Out-Null
